# Edit script generated to apply the diff to before.docx
$d = $word.ActiveDocument

# 1. Title text
$d.Content.Find.Execute("Decoding the Secrets of Immunity", $true, $false, $false, $false, $false, $true, 1, $false, "Chemistry: Unlocking the Secrets of Matter", 2) | Out-Null

# 2. Author name paragraph: "Dr. Sarah Jones" (3 runs) -> "Professor Eleanor Maxwell" (1 run)
$pAuthor = $d.Paragraphs(2)
$rAuthor = $pAuthor.Range
$rAuthor.End = $rAuthor.End - 1
$rAuthor.Text = "Professor Eleanor Maxwell"

# 3. Email address runs
$d.Content.Find.Execute("sjones@healthsciences", $true, $false, $false, $false, $false, $true, 1, $false, "eleanormaxbell6570@gmail", 2) | Out-Null
$d.Content.Find.Execute("edu", $true, $false, $false, $false, $false, $true, 1, $false, "com", 2) | Out-Null

# 4. Body paragraph (paragraph 5): replace entire content (sentences + line breaks) with new Chemistry/Einstein/government text
$pBody = $d.Paragraphs(5)
$rBody = $pBody.Range
$rBody.End = $rBody.End - 1
$rBody.Text = 'Chemistry, the study of matter and its properties, offers a fascinating journey into the world of atoms and molecules, the fundamental building blocks of all substances. It is a science that seeks to understand the interactions between these tiny particles, uncovering the secrets of their behavior and the materials they form. From the air we breathe to the food we eat, chemistry is intricately woven into every aspect of our lives, shaping the world around us in myriad waysOur universe is composed of various elements, the basic substances that cannot be further simplified through chemical means. Chemistry delves into the nature of these elements, exploring their properties, reactivity, and the ways in which they combine to form compounds. By investigating chemical reactions, scientists can manipulate and transform substances, leading to new materials and advancements in various fieldsMoreover, chemistry plays a vital role in our understanding of living organisms. It helps unravel the intricate processes that occur within biological systems, contributing to our knowledge of metabolism, DNA, and the complex interactions that sustain life. Whether it''s the study of biochemical reactions in cells or the development of new drugs, chemistry is essential for advancements in medicine and healthcareEinstein famously said: "The only source of knowledge is experience." Our ability to learn from our experiences has led to a steady accumulation of knowledge over time, helping us understand the world around us better. One area where this is particularly evident is in government, where the study of political systems and public policies can provide valuable lessons for improving governance and decision-making.From the ancient democracies of Greece and Rome to modern representative republics, governments have experimented with various structures and methods of ruling. The study of government helps us understand the strengths and weaknesses of different political systems, allowing us to learn from the successes and failures of the past. This knowledge can inform policy choices and help leaders make more informed decisions, leading to better outcomes for citizens.'

# 5. Summary heading - lastRenderedPageBreak is a rendering artifact; skipped (not settable via object model)

# 6. Summary paragraph (paragraph 7): replace entire content with new Chemistry summary text
$pSummary = $d.Paragraphs(7)
$rSummary = $pSummary.Range
$rSummary.End = $rSummary.End - 1
$rSummary.Text = 'Chemistry, the study of matter and its properties, provides a comprehensive understanding of the world around us. It investigates the behavior of atoms, molecules, and chemical reactions, leading to new materials and advancements in various fields. Chemistry also plays a crucial role in biology, unraveling the intricate processes within living organisms and advancing medicine. Its study offers a valuable lens through which we can understand the world, solve problems, and make informed decisions.'

# 7. Add a new empty paragraph at the very end of the document
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
